$d = $word.ActiveDocument
$d.Content.Find.Execute("Git & Pandoc", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Second 2 测试版本2`rGit & Pandoc", 2)
